$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 207.5
$ws.Range("I28").Value = 207.5
$ws.Range("K28").Value = 207.5
$ws.Range("M28").Value = 277.5
$ws.Range("H32").Value = 1389.2
$ws.Range("I32").Value = 1649
$ws.Range("J32").Value = 999.5
$ws.Range("K32").Value = 1649
$ws.Range("L32").Value = 999.5
$ws.Range("M32").Value = -1323
$ws.Range("N32").Value = -1651.5
$ws.Range("H80").Value = 4352
$ws.Range("I80").Value = 4106
$ws.Range("J80").Value = 4659.5
$ws.Range("K80").Value = 12318
$ws.Range("L80").Value = 13978.5
$ws.Range("M80").Value = -11320
$ws.Range("N80").Value = -15974.5
$ws.Range("H83").Value = 4352
$ws.Range("I83").Value = 4106
$ws.Range("J83").Value = 4659.5
$ws.Range("K83").Value = 36954
$ws.Range("L83").Value = 41935.5
$ws.Range("M83").Value = -31962
$ws.Range("N83").Value = -51919.5
$ws.Range("H94").Value = 1036.6666
$ws.Range("I94").Value = 844
$ws.Range("K94").Value = 844
$ws.Range("M94").Value = -393
$ws.Range("H96").Value = 375.66666
$ws.Range("I96").Value = 235.66667
$ws.Range("K96").Value = 707.00001
$ws.Range("M96").Value = 665.99999
$ws.Range("H98").Value = 2885.7273
$ws.Range("I98").Value = 1174.3
$ws.Range("K98").Value = 1174.3
$ws.Range("M98").Value = 323.7
$ws.Range("H107").Value = 3852.9
$ws.Range("I107").Value = 3316.125
$ws.Range("J107").Value = 6000
$ws.Range("K107").Value = 3316.125
$ws.Range("L107").Value = 6000
$ws.Range("M107").Value = -1396.125
$ws.Range("N107").Value = -9840
$ws.Range("H118").Value = 1259.5
$ws.Range("I118").Value = 1700
$ws.Range("J118").Value = 965.8333
$ws.Range("K118").Value = 5100
$ws.Range("L118").Value = 2897.4999
$ws.Range("M118").Value = -3443
$ws.Range("N118").Value = -6211.4999
$ws.Range("H122").Value = 2885.7273
$ws.Range("I122").Value = 1174.3
$ws.Range("K122").Value = 3522.9
$ws.Range("M122").Value = -1072.9
$ws.Range("H125").Value = 474.75
$ws.Range("I125").Value = 483
$ws.Range("K125").Value = 4347
$ws.Range("M125").Value = -1887
$ws.Range("H132").Value = 1137.2307
$ws.Range("I132").Value = 1177.25
$ws.Range("J132").Value = 657
$ws.Range("K132").Value = 3531.75
$ws.Range("L132").Value = 1971
$ws.Range("M132").Value = -1001.75
$ws.Range("N132").Value = -7031

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 91
$ws.Range("I5").Value = 91
$ws.Range("K5").Value = 91
$ws.Range("M5").Value = 21
$ws.Range("H45").Value = 4099.6
$ws.Range("J45").Value = 1999.5
$ws.Range("L45").Value = 1999.5
$ws.Range("N45").Value = -2753.5
$ws.Range("H97").Value = 3798.125
$ws.Range("I97").Value = 540
$ws.Range("J97").Value = 9228.333000000001
$ws.Range("K97").Value = 540
$ws.Range("L97").Value = 9228.333000000001
$ws.Range("M97").Value = -44
$ws.Range("N97").Value = -10220.333
$ws.Range("H122").Value = 2881.4
$ws.Range("I122").Value = 3151.625
$ws.Range("J122").Value = 1800.5
$ws.Range("K122").Value = 9454.875
$ws.Range("L122").Value = 5401.5
$ws.Range("M122").Value = -7004.875
$ws.Range("N122").Value = -10301.5

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 91
$ws.Range("I4").Value = 91
$ws.Range("K4").Value = 91
$ws.Range("M4").Value = 24
$ws.Range("H20").Value = 3985.6
$ws.Range("I20").Value = 1285.4286
$ws.Range("J20").Value = 6348.25
$ws.Range("K20").Value = 1285.4286
$ws.Range("L20").Value = 6348.25
$ws.Range("M20").Value = -1038.4286
$ws.Range("N20").Value = -6842.25
$ws.Range("H80").Value = 1359.4445
$ws.Range("I80").Value = 881.5
$ws.Range("J80").Value = 1496
$ws.Range("K80").Value = 881.5
$ws.Range("L80").Value = 1496
$ws.Range("M80").Value = 116.5
$ws.Range("N80").Value = -3492
$ws.Range("H83").Value = 1359.4445
$ws.Range("I83").Value = 881.5
$ws.Range("J83").Value = 1496
$ws.Range("K83").Value = 4407.5
$ws.Range("L83").Value = 7480
$ws.Range("M83").Value = 584.5
$ws.Range("N83").Value = -17464
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = $null

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 7500
$ws.Range("I70").Value = 7500
$ws.Range("K70").Value = 7500
$ws.Range("M70").Value = -7185
$ws.Range("H73").Value = 7500
$ws.Range("I73").Value = 7500
$ws.Range("K73").Value = 7500
$ws.Range("M73").Value = -6408
$ws.Range("H99").Value = 8250
$ws.Range("I99").Value = 8250
$ws.Range("K99").Value = 8250
$ws.Range("M99").Value = -6752
$ws.Range("H107").Value = 1461.421
$ws.Range("I107").Value = 1637.1333
$ws.Range("J107").Value = 802.5
$ws.Range("K107").Value = 1637.1333
$ws.Range("L107").Value = 802.5
$ws.Range("M107").Value = 282.8667
$ws.Range("N107").Value = -4642.5
$ws.Range("H122").Value = 5401.8335
$ws.Range("I122").Value = 6128
$ws.Range("K122").Value = 18384
$ws.Range("M122").Value = -15934
$ws.Range("H126").Value = 8250
$ws.Range("I126").Value = 8250
$ws.Range("K126").Value = 24750
$ws.Range("M126").Value = -22280
$ws.Range("H141").Value = 56999
$ws.Range("J141").Value = 56999
$ws.Range("L141").Value = 56999
$ws.Range("N141").Value = -67359

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 4549.857
$ws.Range("I99").Value = 273.5
$ws.Range("K99").Value = 820.5
$ws.Range("M99").Value = 1425.5
$ws.Range("H113").Value = 1276.7142
$ws.Range("I113").Value = 745.25
$ws.Range("J113").Value = 1489.3
$ws.Range("K113").Value = 2235.75
$ws.Range("L113").Value = 4467.9
$ws.Range("M113").Value = -65.75
$ws.Range("N113").Value = -8807.9

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 258.25
$ws.Range("I2").Value = 282
$ws.Range("K2").Value = 282
$ws.Range("M2").Value = -169
$ws.Range("H80").Value = 2250
$ws.Range("J80").Value = 2333.3333
$ws.Range("L80").Value = 2333.3333
$ws.Range("N80").Value = -4329.3333
$ws.Range("H83").Value = 2250
$ws.Range("J83").Value = 2333.3333
$ws.Range("L83").Value = 11666.6665
$ws.Range("N83").Value = -21650.6665
$ws.Range("H102").Value = 6304
$ws.Range("I102").Value = 6304
$ws.Range("K102").Value = 6304
$ws.Range("M102").Value = -4682
$ws.Range("H126").Value = 3956
$ws.Range("I126").Value = 4115.3335
$ws.Range("K126").Value = 12346.0005
$ws.Range("M126").Value = -9876.000499999998

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2676.125
$ws.Range("I7").Value = 2680
$ws.Range("J7").Value = 2649
$ws.Range("K7").Value = 2680
$ws.Range("L7").Value = 2649
$ws.Range("M7").Value = -2568
$ws.Range("N7").Value = -2873
$ws.Range("H40").Value = 6666.6665
$ws.Range("I40").Value = 10000
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 10000
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -9864
$ws.Range("N40").Value = -5272
$ws.Range("H55").Value = 1653.8889
$ws.Range("I55").Value = 2207.8
$ws.Range("J55").Value = 961.5
$ws.Range("K55").Value = 2207.8
$ws.Range("L55").Value = 961.5
$ws.Range("M55").Value = -2034.8
$ws.Range("N55").Value = -1307.5
$ws.Range("H93").Value = 2714.889
$ws.Range("I93").Value = 2714.889
$ws.Range("K93").Value = 2714.889
$ws.Range("M93").Value = -1466.889
$ws.Range("H100").Value = 1745.1818
$ws.Range("I100").Value = 1819.8
$ws.Range("J100").Value = 999
$ws.Range("K100").Value = 1819.8
$ws.Range("L100").Value = 999
$ws.Range("M100").Value = -1278.8
$ws.Range("N100").Value = -2081
$ws.Range("H122").Value = 4424.8
$ws.Range("I122").Value = 4424.8
$ws.Range("K122").Value = 13274.4
$ws.Range("M122").Value = -10824.4
$ws.Range("H126").Value = 2676.125
$ws.Range("I126").Value = 2680
$ws.Range("J126").Value = 2649
$ws.Range("K126").Value = 8040
$ws.Range("L126").Value = 7947
$ws.Range("M126").Value = -5570
$ws.Range("N126").Value = -12887

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1916.6666
$ws.Range("H132").Value = 2145.375
$ws.Range("I132").Value = 1140.0667
$ws.Range("K132").Value = 3420.2001
$ws.Range("M132").Value = -890.2001


Write-Host "Applied all Sophia_Profits market data updates."
